# Apply updates to the "Overview" sheet of the Zagros methanol yearly
# income statement workbook (rial.xlsx).
#
# 1) Update the two "publish date" header labels in row 9 (columns G, H).
# 2) Refresh a batch of computed/reported financial figures for the
#    latest ("1401/12") period column (G) and the H column that holds
#    the newest reported figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Row 9: publish-date header labels ---------------------------------
$ws.Range("G9").Value = "1402-04-14 (9)"
$ws.Range("H9").Value = "1402-04-28 (3)"

# --- Updated figures -----------------------------------------------------
$ws.Range("H12").Value = -145108587
$ws.Range("H13").Value = 69105019

$ws.Range("G14").Value = -30463107
$ws.Range("H14").Value = -44188435

$ws.Range("H16").Value = 5922920

$ws.Range("G17").Value = 15229471
$ws.Range("H17").Value = 30839504

$ws.Range("H19").Value = -479503

$ws.Range("G20").Value = 18966818
$ws.Range("H20").Value = 30360001

$ws.Range("G22").Value = 18966818
$ws.Range("H22").Value = 30360001

$ws.Range("G24").Value = 18966818
$ws.Range("H24").Value = 30360001

$ws.Range("G25").Value = 7903
$ws.Range("H25").Value = 12650

$ws.Range("G27").Value = 7903
$ws.Range("H27").Value = 12650
